# Updating time-varying selex options
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Controls")
$ws.Activate()

# n_sims: 1000 -> 100
$ws.Range("B2").Value = 100

# n_fish_fleets: 1 -> 2
$ws.Range("B6").Value = 2

# Update the active selection shown in the sheet view
$ws.Range("D11").Select()
